$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new columns P and Q, continuing the
# sequence (O1 = 13 -> P1 = 14, Q1 = 15), using the same style as the
# existing header cells (copy format from O1).
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("P1").Value2 = 14
$ws.Range("Q1").Value2 = 15

# For data rows 2-25: swap values in columns I<->K and M<->O, and append
# new columns P and Q with value 2.
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # I
    $kVal = $ws.Cells.Item($r, 11).Value2  # K
    $mVal = $ws.Cells.Item($r, 13).Value2  # M
    $oVal = $ws.Cells.Item($r, 15).Value2  # O

    $ws.Cells.Item($r, 9).Value2 = $kVal   # I = old K
    $ws.Cells.Item($r, 11).Value2 = $iVal  # K = old I
    $ws.Cells.Item($r, 13).Value2 = $oVal  # M = old O
    $ws.Cells.Item($r, 15).Value2 = $mVal  # O = old M

    $ws.Cells.Item($r, 16).Value2 = 2      # P
    $ws.Cells.Item($r, 17).Value2 = 2      # Q
}

$ws.Range("A1").Select() | Out-Null
